# Update p-value columns on each sheet with computed t-test results.
$wb = $excel.ActiveWorkbook

# accuracy sheet
$ws = $wb.Worksheets.Item("accuracy")
$ws.Range("B2").Value = [double]"0.9367569209110885"
$ws.Range("B3").Value = [double]"3.201667760642614e-10"
$ws.Range("B4").Value = [double]"9.581490542972419e-08"
$ws.Range("B5").Value = [double]"0.00399977760565483"
$ws.Range("B6").Value = [double]"0.4545840847298611"
$ws.Range("B10").Value = [double]"0.3415990856234326"
$ws.Range("B12").Value = [double]"0.3978889192397709"
$ws.Range("B15").Value = [double]"0.3107494374775444"

# sensitivity sheet
$ws = $wb.Worksheets.Item("sensitivity")
$ws.Range("B2").Value = [double]"0.2315026934221536"
$ws.Range("B3").Value = [double]"2.09868506100729e-11"
$ws.Range("B4").Value = [double]"3.716408000002603e-09"
$ws.Range("B5").Value = [double]"3.86073535379554e-06"
$ws.Range("B6").Value = [double]"0.07732544779228454"
$ws.Range("B10").Value = [double]"0.0004986954125883854"
$ws.Range("B12").Value = [double]"0.5340768696378337"
$ws.Range("B15").Value = [double]"0.9818449888605499"

# specificity sheet
$ws = $wb.Worksheets.Item("specificity")
$ws.Range("B2").Value = [double]"0.376069952301525"
$ws.Range("B3").Value = [double]"0.02215798961689382"
$ws.Range("B4").Value = [double]"0.02468624933218203"
$ws.Range("B5").Value = [double]"0.1291322430291489"
$ws.Range("B6").Value = [double]"0.003246495043908954"
$ws.Range("B10").Value = [double]"6.098849253641642e-05"
$ws.Range("B12").Value = [double]"0.1625510095051058"
$ws.Range("B15").Value = [double]"0.007682017017812855"

# time sheet - replace all p-values with the new computed value
$ws = $wb.Worksheets.Item("time")
$newTimeValue = [double]"2.009597066272628e-33"
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 2).Value = $newTimeValue
}
